$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Change 1: "...Asia, mucho de los cuales..." -> "...Asia, muchos de
# los cuales..." (typo fix) and the _GoBack bookmark is relocated here
# (right after the newly typed "s"), splitting the original run into
# three runs: "...mucho" | "s" | " de ".
# ------------------------------------------------------------------

# Place a transient "protective" bookmark at the boundary right before
# "los cuales..." so that the engine's run-coalescing (which normally
# merges every adjacent, identically-formatted run in the paragraph
# whenever the paragraph is edited) stops there and leaves the later
# runs of the paragraph untouched.
$rProtect = $d.Content
$rProtect.Find.Execute("los cuales son una herencia", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$protectPos = $rProtect.Start
$protectRange = $d.Range($protectPos, $protectPos)
$d.Bookmarks.Add("ZZPROTECT1", $protectRange)

# Locate "mucho de" and compute the insertion point right after "mucho"
# (5 characters) where the new "s" must be typed.
$rFind = $d.Content
$rFind.Find.Execute("mucho de", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$sPos = $rFind.Start + 5

# Type the "s".
$insRange = $d.Range($sPos, $sPos)
$insRange.InsertAfter("s")

# Toggling formatting on/off forces the newly typed character to stay
# in its own run instead of being re-absorbed into the previous run.
$newCharRange = $d.Range($sPos, $sPos + 1)
$newCharRange.Bold = 1
$newCharRange.Bold = 0

# Move the _GoBack bookmark to sit right after the inserted "s" (this
# is where Word leaves it after the most recent edit).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$goBackPos = $sPos + 1
$goBackRange = $d.Range($goBackPos, $goBackPos)
$d.Bookmarks.Add("_GoBack", $goBackRange)

# Remove the transient protective bookmark, its job is done.
$d.Bookmarks.Item("ZZPROTECT1").Delete()

# ------------------------------------------------------------------
# Change 2: the old _GoBack location ("de manera" | " cronológica.")
# no longer has a bookmark in between, so the two runs collapse back
# into a single "de manera cronológica." run.
# ------------------------------------------------------------------
$rMerge = $d.Content
$rMerge.Find.Execute(" cronológica.", $true, $false, $false, $false, $false, $true, 1, $false, " cronológica.", 2)
